$wb = $excel.ActiveWorkbook

# Sheet 1: "Potencia Acumulada - SIN (MW)"
# - E1 label was a stray numeric value; should read "2050" (text, matching the
#   other year-header cells in row 1)
# - drop the trailing "Total" row (row 13)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").Value = "'2050"
$ws1.Rows.Item(13).Delete()

# Sheet 2: "Geracao Periodo Medio (MWMed)"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").Value = "'2050"
$ws2.Rows.Item(13).Delete()

# Sheet 3: "Atendimento a Ponta(MW)"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").Value = "'2050"
$ws3.Rows.Item(13).Delete()

# Sheet 4: "Potencia Incremental - SIN(MW)" — header row uses ranges (e.g.
# "2031-2040"), so the fixed E1 label is "2041-2050"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" — same stray-value fix, no Total row here
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "'2050"

# Sheet 6: "Custo Total (bilhões de R$)" — no E1 label on this sheet, just
# drop the trailing "Total" row (row 4)
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
